$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    # All data cells in this sheet are stored as text, even the
    # numeric-looking ones (Quantity/Price/TotalIncome columns).
    # A bare numeric literal would be auto-coerced by Excel into a
    # Number cell, so numeric-looking values get a leading apostrophe
    # to force text entry (the apostrophe itself is not stored).
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Cells.Item($row, $col).Value = "'" + $val
    } else {
        $ws.Cells.Item($row, $col).Value = $val
    }
}

# --- Existing rows: update changed cells ---

# Row 2: AIKO / ARONChair -> IKEAChair, 70 -> 60, 700 -> 600
Set-TextCell 2 2 "IKEAChair"
Set-TextCell 2 4 "60"
Set-TextCell 2 5 "600"

# Row 3: Telerik Academy / Pesho -> IKEABed, 15.5 -> 120, 806 -> 6240
Set-TextCell 3 2 "IKEABed"
Set-TextCell 3 4 "120"
Set-TextCell 3 5 "6240"

# Row 4: Boyana Film Studios / Pesho -> ARONMirror, 20 -> 18, 480 -> 432
Set-TextCell 4 2 "ARONMirror"
Set-TextCell 4 4 "18"
Set-TextCell 4 5 "432"

# Row 5: San Benedetto / Pesho -> VidenovChair (qty/price/total unchanged)
Set-TextCell 5 2 "VidenovChair"

# Row 6: Mladost Estate / ARONMirror -> ARONChair, 18 -> 70, 72 -> 280
Set-TextCell 6 2 "ARONChair"
Set-TextCell 6 4 "70"
Set-TextCell 6 5 "280"

# --- New rows 7-11 ---

# Row 7: Sunset Security / IKEAMirror / 14 / 20 / 280
Set-TextCell 7 1 "Sunset Security"
Set-TextCell 7 2 "IKEAMirror"
Set-TextCell 7 3 "14"
Set-TextCell 7 4 "20"
Set-TextCell 7 5 "280"

# Row 8: West Bank / ARONChair / 20 / 70 / 1400
Set-TextCell 8 1 "West Bank"
Set-TextCell 8 2 "ARONChair"
Set-TextCell 8 3 "20"
Set-TextCell 8 4 "70"
Set-TextCell 8 5 "1400"

# Row 9: BILLA / ARONBed / 17 / 150 / 2550
Set-TextCell 9 1 "BILLA"
Set-TextCell 9 2 "ARONBed"
Set-TextCell 9 3 "17"
Set-TextCell 9 4 "150"
Set-TextCell 9 5 "2550"

# Row 10: Null Industries / IKEABed / 21 / 120 / 2520
Set-TextCell 10 1 "Null Industries"
Set-TextCell 10 2 "IKEABed"
Set-TextCell 10 3 "21"
Set-TextCell 10 4 "120"
Set-TextCell 10 5 "2520"

# Row 11: VS Incorporated / IKEAMirror / 7 / 20 / 140
Set-TextCell 11 1 "VS Incorporated"
Set-TextCell 11 2 "IKEAMirror"
Set-TextCell 11 3 "7"
Set-TextCell 11 4 "20"
Set-TextCell 11 5 "140"

# --- Selection moves to G6 ---
[void]$ws.Range("G6").Select()
